$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(36, 2.123996793003155, -1.341851536611955, 0.2826872811415794, 0.4676428022577083, 0.3615278921897729),
    @(37, 2.999549018034927, -1.305966639919423, 0.2789175382363926, 0.4622694510566538, 0.3606330597665377),
    @(38, 2.006392134354136, -1.342358922876254, 0.2780857581067647, 0.4626607278200771, 0.3603571426381924),
    @(39, 3.240696515968457, -1.295640757642037, 0.2802419799738735, 0.4654042137349904, 0.3627833696019208),
    @(40, 3.730790296199632, -1.274821990365726, 0.2861881056395761, 0.4688292432995688, 0.3644847951351571),
    @(41, 7.168579958009738, -1.124017827048016, 0.2972256193218946, 0.4943963225100994, 0.3866122479142569),
    @(42, 6.238961826588797, -1.164018083528414, 0.293609252693852,  0.4761356968837653, 0.3778501832234241),
    @(43, 0.5746675060904044, -1.410359791612187, 0.2727731295019905, 0.4610450423698221, 0.3518740037090754),
    @(44, 11.18770091948372, -0.9351897937297551, 0.3124566566827671, 0.5171976521000549, 0.4083378066995465),
    @(45, 3.454749111365462, -1.286853781347119, 0.2846692163941658, 0.4718034380622428, 0.3670012483452475),
    @(46, 10.74914869260767, -0.967016390682338, 0.3101819694086514, 0.5107950315223438, 0.4046876801675777)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}

# Copy the style from A10 (existing bordered/bold/centered style) to the new A11, A12 cells
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A12").PasteSpecial(-4122) | Out-Null
